$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1387.2982
$ws.Cells.Item(15, 9).Value = 1387.2982
$ws.Cells.Item(15, 11).Value = 4161.8946
$ws.Cells.Item(15, 13).Value = -3992.8946
$ws.Cells.Item(28, 8).Value = 941.2727
$ws.Cells.Item(28, 9).Value = 483.77777
$ws.Cells.Item(28, 10).Value = 3000
$ws.Cells.Item(28, 11).Value = 483.77777
$ws.Cells.Item(28, 12).Value = 3000
$ws.Cells.Item(28, 13).Value = 1.222230000000025
$ws.Cells.Item(28, 14).Value = -3970
$ws.Cells.Item(113, 8).Value = 60376.35
$ws.Cells.Item(113, 9).Value = 126149.75
$ws.Cells.Item(113, 10).Value = 1911.1111
$ws.Cells.Item(113, 11).Value = 126149.75
$ws.Cells.Item(113, 12).Value = 1911.1111
$ws.Cells.Item(113, 13).Value = -122895.75
$ws.Cells.Item(113, 14).Value = -8419.1111
$ws.Cells.Item(129, 8).Value = 857.5952
$ws.Cells.Item(129, 9).Value = 579.2
$ws.Cells.Item(129, 10).Value = 1012.2593
$ws.Cells.Item(129, 11).Value = 1737.6
$ws.Cells.Item(129, 12).Value = 3036.7779
$ws.Cells.Item(129, 13).Value = 3262.4
$ws.Cells.Item(129, 14).Value = -13036.7779
$ws.Cells.Item(132, 8).Value = 4060.5789
$ws.Cells.Item(132, 9).Value = 4603.4062
$ws.Cells.Item(132, 10).Value = 1165.5
$ws.Cells.Item(132, 11).Value = 13810.2186
$ws.Cells.Item(132, 12).Value = 3496.5
$ws.Cells.Item(132, 13).Value = -11280.2186
$ws.Cells.Item(132, 14).Value = -8556.5
$ws.Cells.Item(135, 8).Value = 874.25
$ws.Cells.Item(135, 9).Value = 960.86664
$ws.Cells.Item(135, 11).Value = 8647.79976
$ws.Cells.Item(135, 13).Value = -6112.79976
$ws.Cells.Item(137, 8).Value = 1377.0962
$ws.Cells.Item(137, 9).Value = 995.7692
$ws.Cells.Item(137, 10).Value = 1504.2051
$ws.Cells.Item(137, 11).Value = 2987.3076
$ws.Cells.Item(137, 12).Value = 4512.615299999999
$ws.Cells.Item(137, 13).Value = -437.3076000000001
$ws.Cells.Item(137, 14).Value = -9612.615299999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2168.375
$ws.Cells.Item(61, 9).Value = 1202.7
$ws.Cells.Item(61, 10).Value = 2858.1428
$ws.Cells.Item(61, 11).Value = 1202.7
$ws.Cells.Item(61, 12).Value = 2858.1428
$ws.Cells.Item(61, 13).Value = -990.7
$ws.Cells.Item(61, 14).Value = -3282.1428
$ws.Cells.Item(63, 8).Value = 2999.923
$ws.Cells.Item(63, 9).Value = 2257
$ws.Cells.Item(63, 10).Value = 3866.6667
$ws.Cells.Item(63, 11).Value = 2257
$ws.Cells.Item(63, 12).Value = 3866.6667
$ws.Cells.Item(63, 13).Value = -1571
$ws.Cells.Item(63, 14).Value = -5238.6667
$ws.Cells.Item(66, 8).Value = 2999.923
$ws.Cells.Item(66, 9).Value = 2257
$ws.Cells.Item(66, 10).Value = 3866.6667
$ws.Cells.Item(66, 11).Value = 11285
$ws.Cells.Item(66, 12).Value = 19333.3335
$ws.Cells.Item(66, 13).Value = -7853
$ws.Cells.Item(66, 14).Value = -26197.3335
$ws.Cells.Item(74, 8).Value = 1313.8085
$ws.Cells.Item(74, 9).Value = 824.8333
$ws.Cells.Item(74, 10).Value = 2176.7058
$ws.Cells.Item(74, 11).Value = 824.8333
$ws.Cells.Item(74, 12).Value = 2176.7058
$ws.Cells.Item(74, 13).Value = 49.16669999999999
$ws.Cells.Item(74, 14).Value = -3924.7058
$ws.Cells.Item(77, 8).Value = 1313.8085
$ws.Cells.Item(77, 9).Value = 824.8333
$ws.Cells.Item(77, 10).Value = 2176.7058
$ws.Cells.Item(77, 11).Value = 4124.1665
$ws.Cells.Item(77, 12).Value = 10883.529
$ws.Cells.Item(77, 13).Value = 243.8334999999997
$ws.Cells.Item(77, 14).Value = -19619.529
$ws.Cells.Item(88, 8).Value = 13033
$ws.Cells.Item(88, 9).Value = 1899
$ws.Cells.Item(88, 11).Value = 1899
$ws.Cells.Item(88, 13).Value = -1493
$ws.Cells.Item(91, 8).Value = 13033
$ws.Cells.Item(91, 9).Value = 1899
$ws.Cells.Item(91, 11).Value = 1899
$ws.Cells.Item(91, 13).Value = -495
$ws.Cells.Item(132, 8).Value = 2539.9805
$ws.Cells.Item(132, 9).Value = 2448.7708
$ws.Cells.Item(132, 10).Value = 3999.3333
$ws.Cells.Item(132, 11).Value = 7346.312399999999
$ws.Cells.Item(132, 12).Value = 11997.9999
$ws.Cells.Item(132, 13).Value = -4816.312399999999
$ws.Cells.Item(132, 14).Value = -17057.9999
$ws.Cells.Item(136, 8).Value = 2168.375
$ws.Cells.Item(136, 9).Value = 1202.7
$ws.Cells.Item(136, 10).Value = 2858.1428
$ws.Cells.Item(136, 11).Value = 3608.1
$ws.Cells.Item(136, 12).Value = 8574.428400000001
$ws.Cells.Item(136, 13).Value = -1058.1
$ws.Cells.Item(136, 14).Value = -13674.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1512.1351
$ws.Cells.Item(134, 9).Value = 1470.8055
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 4412.416499999999
$ws.Cells.Item(134, 12).Value = 9000
$ws.Cells.Item(134, 13).Value = -1877.416499999999
$ws.Cells.Item(134, 14).Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2664
$ws.Cells.Item(132, 9).Value = 1604.4
$ws.Cells.Item(132, 10).Value = 3988.5
$ws.Cells.Item(132, 11).Value = 4813.200000000001
$ws.Cells.Item(132, 12).Value = 11965.5
$ws.Cells.Item(132, 13).Value = -2283.200000000001
$ws.Cells.Item(132, 14).Value = -17025.5
$ws.Cells.Item(134, 8).Value = 1596.8
$ws.Cells.Item(134, 9).Value = 1436.5
$ws.Cells.Item(134, 11).Value = 4309.5
$ws.Cells.Item(134, 13).Value = -1774.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1175.6666
$ws.Cells.Item(34, 9).Value = 301
$ws.Cells.Item(34, 11).Value = 903
$ws.Cells.Item(34, 13).Value = -819
$ws.Cells.Item(68, 8).Value = 18738.88
$ws.Cells.Item(68, 9).Value = 1115.2667
$ws.Cells.Item(68, 10).Value = 24746.932
$ws.Cells.Item(68, 11).Value = 3345.800099999999
$ws.Cells.Item(68, 12).Value = 74240.796
$ws.Cells.Item(68, 13).Value = -2534.800099999999
$ws.Cells.Item(68, 14).Value = -75862.796
$ws.Cells.Item(71, 8).Value = 18738.88
$ws.Cells.Item(71, 9).Value = 1115.2667
$ws.Cells.Item(71, 10).Value = 24746.932
$ws.Cells.Item(71, 11).Value = 10037.4003
$ws.Cells.Item(71, 12).Value = 222722.388
$ws.Cells.Item(71, 13).Value = -5981.400299999999
$ws.Cells.Item(71, 14).Value = -230834.388
$ws.Cells.Item(131, 8).Value = 895.26154
$ws.Cells.Item(131, 10).Value = 944.9138
$ws.Cells.Item(131, 12).Value = 2834.7414
$ws.Cells.Item(131, 14).Value = -12914.7414

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1432
$ws.Cells.Item(113, 9).Value = 1300
$ws.Cells.Item(113, 10).Value = 1469.7142
$ws.Cells.Item(113, 11).Value = 1300
$ws.Cells.Item(113, 12).Value = 1469.7142
$ws.Cells.Item(113, 13).Value = 870
$ws.Cells.Item(113, 14).Value = -5809.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 183334.67
$ws.Cells.Item(14, 10).Value = 75002
$ws.Cells.Item(14, 12).Value = 75002
$ws.Cells.Item(14, 14).Value = -75346
$ws.Cells.Item(16, 8).Value = 63150.625
$ws.Cells.Item(16, 9).Value = 71900.71000000001
$ws.Cells.Item(16, 10).Value = 1900
$ws.Cells.Item(16, 11).Value = 71900.71000000001
$ws.Cells.Item(16, 12).Value = 1900
$ws.Cells.Item(16, 13).Value = -71730.71000000001
$ws.Cells.Item(16, 14).Value = -2240
$ws.Cells.Item(132, 8).Value = 5645.409
$ws.Cells.Item(132, 9).Value = 6800.0835
$ws.Cells.Item(132, 10).Value = 4259.8
$ws.Cells.Item(132, 11).Value = 20400.2505
$ws.Cells.Item(132, 12).Value = 12779.4
$ws.Cells.Item(132, 13).Value = -17870.2505
$ws.Cells.Item(132, 14).Value = -17839.4
$ws.Cells.Item(136, 8).Value = 1101.2333
$ws.Cells.Item(136, 9).Value = 888.8333
$ws.Cells.Item(136, 11).Value = 2666.4999
$ws.Cells.Item(136, 13).Value = -116.4998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2059
$ws.Cells.Item(132, 9).Value = 1575.3572
$ws.Cells.Item(132, 10).Value = 2905.375
$ws.Cells.Item(132, 11).Value = 4726.071599999999
$ws.Cells.Item(132, 12).Value = 8716.125
$ws.Cells.Item(132, 13).Value = -2196.071599999999
$ws.Cells.Item(132, 14).Value = -13776.125
$ws.Cells.Item(136, 8).Value = 1270.1904
$ws.Cells.Item(136, 9).Value = 583.9286
$ws.Cells.Item(136, 11).Value = 1751.7858
$ws.Cells.Item(136, 13).Value = 798.2142000000001
